# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (with per-fund holding detail) between the
# existing "2021-Q4" sheet and the "总计" (total) summary sheet, and updates
# the "总计" sheet with a new leading row summarising the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right after "2021-Q4".
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# NOTE: fetch "总计" only now (by name), since inserting a new sheet shifts
# the positional index every other sheet handle would otherwise resolve to.
$wsTotal = $wb.Worksheets.Item("总计")

# Pick up the same look & feel as the "2021-Q4" sheet: bold/bordered header
# row, and bold/bordered index column (column A).
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A14").PasteSpecial(-4122)

# Header row
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Columns B, C, D, E, F, G hold text (even when the text looks numeric, e.g.
# "10.67"), so force text formatting before writing those values. Column A
# (row index) and column H (rank) are genuine numbers.
$wsQ1.Range("B2:G14").NumberFormat = "@"

$rows = @(
    @{ A = 0;  B = "001070"; C = "建信信息产业股票";                     D = "10.67"; E = "81.80"; F = "6.57"; G = "0.7010"; H = 2 },
    @{ A = 1;  B = "008962"; C = "建信科技创新混合A";                    D = "4.34";  E = "84.82"; F = "6.62"; G = "0.2873"; H = 2 },
    @{ A = 2;  B = "000308"; C = "建信创新中国混合";                     D = "3.11";  E = "84.50"; F = "6.50"; G = "0.2022"; H = 2 },
    @{ A = 3;  B = "010994"; C = "博时创新经济混合A";                    D = "3.89";  E = "89.89"; F = "4.76"; G = "0.1852"; H = 6 },
    @{ A = 4;  B = "002378"; C = "建信弘利灵活配置混合";                  D = "1.03";  E = "89.57"; F = "3.84"; G = "0.0396"; H = 6 },
    @{ A = 5;  B = "002281"; C = "建信裕利灵活配置混合";                  D = "1.10";  E = "88.94"; F = "3.50"; G = "0.0385"; H = 7 },
    @{ A = 6;  B = "970020"; C = "信达价值精选一年持有期灵活配置混合A";    D = "0.64";  E = "56.02"; F = "5.62"; G = "0.0360"; H = 3 },
    @{ A = 7;  B = "970021"; C = "信达价值精选一年持有期灵活配置混合B";    D = "0.53";  E = "56.02"; F = "5.62"; G = "0.0298"; H = 3 },
    @{ A = 8;  B = "008963"; C = "建信科技创新混合C";                    D = "0.26";  E = "84.82"; F = "6.62"; G = "0.0172"; H = 2 },
    @{ A = 9;  B = "010995"; C = "博时创新经济混合C";                    D = "0.35";  E = "89.89"; F = "4.76"; G = "0.0167"; H = 6 },
    @{ A = 10; B = "519951"; C = "长信利泰灵活配置混合A";                 D = "0.07";  E = "25.18"; F = "0.70"; G = "0.0005"; H = 7 },
    @{ A = 11; B = "008071"; C = "长信利泰灵活配置混合E";                 D = "0.01";  E = "25.18"; F = "0.70"; G = "0.0001"; H = 7 },
    @{ A = 12; B = "007863"; C = "长信利泰灵活配置混合C";                 D = "0.00";  E = "25.18"; F = "0.70"; G = "SKIP";    H = 7 }
)

$r = 2
foreach ($row in $rows) {
    $wsQ1.Cells.Item($r, 1).Value = $row.A
    $wsQ1.Cells.Item($r, 2).Value = $row.B
    $wsQ1.Cells.Item($r, 3).Value = $row.C
    $wsQ1.Cells.Item($r, 4).Value = $row.D
    $wsQ1.Cells.Item($r, 5).Value = $row.E
    $wsQ1.Cells.Item($r, 6).Value = $row.F
    if ($row.G -ne "SKIP") {
        $wsQ1.Cells.Item($r, 7).Value = $row.G
    }
    $wsQ1.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Last data row's holding-value (G14) is a genuine 0, stored as a real number
# (not text like the other G-column cells above it).
$wsQ1.Cells.Item(14, 7).NumberFormat = "General"
$wsQ1.Cells.Item(14, 7).Value = 0

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row for "2022-Q1" above the
#    existing "2021-Q4" row, so the latest quarter sits on top.
# ---------------------------------------------------------------------
$wsTotal.Range("A2:A2").EntireRow.Insert()

# Row-insert borrows formatting from neighbouring rows inconsistently (the
# index column ends up unstyled while the data columns end up bold), so
# reset it explicitly: data cells (B:D) get plain/default formatting, and
# the index cell (A) gets the same bold+bordered look as the row below it.
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 13
$wsTotal.Cells.Item(2, 4).Value = 1.55

$wsTotal.Cells.Item(3, 1).Value = 1
